# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking "Price" values must stay text (as in the source data), so
# those assignments are apostrophe-prefixed to force text entry and then the
# cell style is reset to "Normal" to avoid leaving a stray quote-prefix style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.211.78"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "2.482.63"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'563.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "'162.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.49%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "2.481.00"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").Value = "'0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.61%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "'0.350"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.81%  "
$ws.Range("D13").Value = "'4.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "69.071.26"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "'24.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.05%  "
$ws.Range("D18").Value = "2.479.48"
$ws.Range("E18").Value = "  -3.12%  "
$ws.Range("D19").Value = "'11.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.88%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.70%  "
$ws.Range("D21").Value = "'342.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.66%  "
$ws.Range("D22").Value = "'3.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "'1.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.99%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'69.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "'3.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.17%  "
$ws.Range("D27").Value = "2.609.51"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").Value = "'8.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "0.0₃0862"
$ws.Range("E30").Value = "  -6.33%  "
$ws.Range("D31").Value = "'7.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.42%  "
$ws.Range("D32").Value = "'438.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.84%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.15%  "
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").Value = "'154.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").Value = "'0.112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.12%  "
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").Value = "'18.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("D42").Value = "'4.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.58%  "
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.97%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'1.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.09%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'138.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0724"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.569"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0917"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "
